$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting/style from the last existing row (712) down through the new rows (713:734)
$ws.Range("A712:C712").Copy($ws.Range("A713:C734"))

# Populate the new rows with their data
$ws.Cells.Item(713, 1).Value = 44175
$ws.Cells.Item(713, 2).Value = 29090.04
$ws.Cells.Item(713, 3).Value = 30182.71

$ws.Cells.Item(714, 1).Value = 44176
$ws.Cells.Item(714, 2).Value = 29089.1
$ws.Cells.Item(714, 3).Value = 30184.92

$ws.Cells.Item(715, 1).Value = 44177
$ws.Cells.Item(715, 2).Value = 29088.16
$ws.Cells.Item(715, 3).Value = 30187.14

$ws.Cells.Item(716, 1).Value = 44178
$ws.Cells.Item(716, 2).Value = 29087.22
$ws.Cells.Item(716, 3).Value = 30189.360000000001

$ws.Cells.Item(717, 1).Value = 44179
$ws.Cells.Item(717, 2).Value = 29086.29
$ws.Cells.Item(717, 3).Value = 30191.57

$ws.Cells.Item(718, 1).Value = 44180
$ws.Cells.Item(718, 2).Value = 29085.35
$ws.Cells.Item(718, 3).Value = 30193.79

$ws.Cells.Item(719, 1).Value = 44181
$ws.Cells.Item(719, 2).Value = 29084.41
$ws.Cells.Item(719, 3).Value = 30196.01

$ws.Cells.Item(720, 1).Value = 44182
$ws.Cells.Item(720, 2).Value = 29083.47
$ws.Cells.Item(720, 3).Value = 30198.23

$ws.Cells.Item(721, 1).Value = 44183
$ws.Cells.Item(721, 2).Value = 29082.53
$ws.Cells.Item(721, 3).Value = 30200.44

$ws.Cells.Item(722, 1).Value = 44184
$ws.Cells.Item(722, 2).Value = 29081.59
$ws.Cells.Item(722, 3).Value = 30202.66

$ws.Cells.Item(723, 1).Value = 44185
$ws.Cells.Item(723, 2).Value = 29080.65
$ws.Cells.Item(723, 3).Value = 30204.880000000001

$ws.Cells.Item(724, 1).Value = 44186
$ws.Cells.Item(724, 2).Value = 29079.72
$ws.Cells.Item(724, 3).Value = 30207.1

$ws.Cells.Item(725, 1).Value = 44187
$ws.Cells.Item(725, 2).Value = 29078.78
$ws.Cells.Item(725, 3).Value = 30209.31

$ws.Cells.Item(726, 1).Value = 44188
$ws.Cells.Item(726, 2).Value = 29077.84
$ws.Cells.Item(726, 3).Value = 30211.53

$ws.Cells.Item(727, 1).Value = 44189
$ws.Cells.Item(727, 2).Value = 29076.9
$ws.Cells.Item(727, 3).Value = 30213.75

$ws.Cells.Item(728, 1).Value = 44190
$ws.Cells.Item(728, 2).Value = 29075.96
$ws.Cells.Item(728, 3).Value = 30215.97

$ws.Cells.Item(729, 1).Value = 44191
$ws.Cells.Item(729, 2).Value = 29075.02
$ws.Cells.Item(729, 3).Value = 30218.19

$ws.Cells.Item(730, 1).Value = 44192
$ws.Cells.Item(730, 2).Value = 29074.080000000002
$ws.Cells.Item(730, 3).Value = 30220.41

$ws.Cells.Item(731, 1).Value = 44193
$ws.Cells.Item(731, 2).Value = 29073.15
$ws.Cells.Item(731, 3).Value = 30222.63

$ws.Cells.Item(732, 1).Value = 44194
$ws.Cells.Item(732, 2).Value = 29072.21
$ws.Cells.Item(732, 3).Value = 30224.85

$ws.Cells.Item(733, 1).Value = 44195
$ws.Cells.Item(733, 2).Value = 29071.27
$ws.Cells.Item(733, 3).Value = 30227.07

$ws.Cells.Item(734, 1).Value = 44196
$ws.Cells.Item(734, 2).Value = 29070.33
$ws.Cells.Item(734, 3).Value = 30229.29

# Update the defined name range to cover the new data (A1:C734)
$wb.Names.Item(1).RefersTo = "=UF_IVP_DIARIO!`$A`$1:`$C`$734"

# Reposition frozen-pane view / active selection to mirror the refreshed data range
$ws.Range("B734").Select()
